$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 161; this shifts existing rows 161..247 down to 162..248
$ws.Rows.Item(161).Insert()

# Populate the newly inserted row 161 with the new record's data.
$ws.Range("A161").Value = 8
$ws.Range("B161").Value = "Terminal La Palmera de La Serena"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 45089
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 100112001
$ws.Range("G161").Value = "Berenjena"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Primera"
$ws.Range("J161").Value = 200
$ws.Range("K161").Value = 8500
$ws.Range("L161").Value = 9000
$ws.Range("M161").Value = 8750
$ws.Range("N161").Value = "$/caja 50 unidades"
$ws.Range("O161").Value = "Región de Arica y Parinacota"
$ws.Range("P161").Value = 175
$ws.Range("Q161").Value = 50
$ws.Range("R161").Value = "Hortaliza"
